$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Project ID for the second officer registration (row 3) was corrected
# from 1 to 3.
$ws.Range("C3").Value = 3

# Leave the selection on the cell that was just edited.
$ws.Range("C3").Select()
